# "Added Scrum report, Sprint backlog"
#
# The Sprint 1 backlog sheet tracked which sprint each story belonged to in
# column C ("Sprint"). The four highest-priority stories (rows 8-11, IDs
# 1, 2, 3 and 5) were originally scoped into Sprint 1 but are now moved
# into Sprint 2 as part of building out the Sprint backlog / scrum report.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 1")

$ws.Range("C8").Value  = 2
$ws.Range("C9").Value  = 2
$ws.Range("C10").Value = 2
$ws.Range("C11").Value = 2

# Reflect the user's current focus on the backlog view: scrolled down a bit
# and zoomed in on the newly-updated rows.
$ws.Range("B10").Select()
$excel.ActiveWindow.Zoom = 150
